# Incorporo nuevos datos hasta diciembre de 2025
# Adds 5 new monthly columns (2025M08..2025M12) after the existing last
# data column (JX = 2025M07) to the "tabla-50913" sheet, for every CCAA row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clone the formatting of the current last data column (JX, col 284)
#        onto the five new columns (JY:KC, cols 285-289) so the new cells
#        pick up the same header style (s=2) and numeric style (s=4).
$ws.Range("JX1:JX20").Copy() | Out-Null
$ws.Range("JY1:KC20").PasteSpecial(-4122) | Out-Null

# --- 2. New header labels (row 1) for the five new month columns.
$newHeaders = @{
    285 = "2025M08"
    286 = "2025M09"
    287 = "2025M10"
    288 = "2025M11"
    289 = "2025M12"
}
foreach ($col in $newHeaders.Keys) {
    $ws.Cells.Item(1, $col).Value = $newHeaders[$col]
}

# --- 3. New IPC values for each CCAA row (2..20), columns JY,JZ,KA,KB,KC.
$newMonthData = @{
    2  = @(119.309, 119.106, 120.163, 120.139, 120.524)
    3  = @(118.912, 118.557, 119.472, 119.711, 120.016)
    4  = @(119.688, 118.655, 119.205, 119.375, 119.791)
    5  = @(120.41,  119.683, 120.201, 120.03,  120.332)
    6  = @(118.537, 117.759, 118.784, 119,     119.798)
    7  = @(118.985, 118.101, 119.125, 119.524, 119.991)
    8  = @(119.437, 119.033, 119.889, 120.411, 120.713)
    9  = @(119.932, 119.641, 120.695, 120.927, 121.476)
    10 = @(118.248, 117.655, 118.162, 118.46,  118.829)
    11 = @(119.458, 118.996, 119.767, 120.04,  120.291)
    12 = @(118.832, 118.385, 119.405, 119.926, 120.159)
    13 = @(119.519, 119.026, 119.798, 120.286, 120.469)
    14 = @(117.28,  117.639, 118.567, 118.712, 119.383)
    15 = @(118.498, 118.091, 118.663, 119.321, 119.672)
    16 = @(119.283, 118.45,  119.631, 119.862, 120.082)
    17 = @(119.35,  118.996, 119.698, 119.941, 120.22)
    18 = @(118.83,  118.37,  119.12,  119.384, 120.117)
    19 = @(119.115, 118.815, 119.088, 119.237, 120.023)
    20 = @(119.754, 118.883, 120.284, 120.468, 121.029)
}

$newCols = @(285, 286, 287, 288, 289)

foreach ($row in $newMonthData.Keys) {
    $vals = $newMonthData[$row]
    for ($i = 0; $i -lt $newCols.Length; $i++) {
        $ws.Cells.Item($row, $newCols[$i]).Value = $vals[$i]
    }
}

# --- 4. Match the author's final view state (scrolled right, cell KB27
#        selected) as closely as the object model allows.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 278
$ws.Range("KB27").Select() | Out-Null
